# Weekly update: insert one new price-record row for "Pepino ensalada"
# (Macroferia Regional de Talca) ahead of the existing row 193, shifting
# the rest of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 193; everything currently at 193..211 moves to 194..212.
$ws.Rows.Item(193).Insert()

# Populate the newly-inserted row 193 with the latest weekly record.
$ws.Cells.Item(193, 1).Value = 5
$ws.Cells.Item(193, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(193, 3).Value = "Maule"
$ws.Cells.Item(193, 4).Value = 44449
$ws.Cells.Item(193, 5).Value = 7
$ws.Cells.Item(193, 6).Value = 100112043
$ws.Cells.Item(193, 7).Value = "Pepino ensalada"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Primera"
$ws.Cells.Item(193, 10).Value = 300
$ws.Cells.Item(193, 11).Value = 16000
$ws.Cells.Item(193, 12).Value = 16000
$ws.Cells.Item(193, 13).Value = 16000
$ws.Cells.Item(193, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(193, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(193, 16).Value = 267
$ws.Cells.Item(193, 17).Value = 60
$ws.Cells.Item(193, 18).Value = "Hortaliza"
